$wb = $excel.ActiveWorkbook

# --- "metric overrides" sheet: update the stored selection ---
$wsMetric = $wb.Worksheets.Item("metric overrides")
$wsMetric.Activate()
$wsMetric.Range("D19").Select()

# --- "production" sheet: add a new production item / sub item type column ---
$ws = $wb.Worksheets.Item("production")
$ws.Activate()

# Insert a new row for the new mortgage portfolio item, right after the first data row
$ws.Rows("4:4").Insert()

$ws.Range("A4").Value = $ws.Range("A3").Value2
$ws.Range("B4").Value = $ws.Range("B3").Text
$ws.Range("C4").Value = $ws.Range("C3").Value2
$ws.Range("D4").Value = $ws.Range("D3").Value2

# New "Sub Item Type" column so new labels can be introduced per production item
$ws.Range("E2").Value = "Sub Item Type"
$ws.Range("E4").Value = "New mortgage portfolio"

$ws.Range("E3").Select()

Write-Host "done"
